# Reset the "Returned to White House" counts to 0 across each category.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("B51").Value = 0
